# Generate Report for Handback
# A new handoff/handback cycle is recorded for the "075d7a73-..." file in
# both the zh-cn and de-de status tables: its "Correspond Handoff Datetime"
# (column E) and "Correspond Handback DateTime" (column H) on row 2 are
# refreshed to new timestamps, while the "1bbacf06-..." file's row (row 3)
# keeps its existing timestamps unchanged.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 10:37:17"
$wsZhCn.Range("H2").Value = "2016-03-19 10:37:36"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 10:37:20"
$wsDeDe.Range("H2").Value = "2016-03-19 10:37:41"
